$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "TheFrozenGeneral"
$ws.Range("C39").Value = "The Frozen General"
$ws.Range("D39").Value = 2
$ws.Range("E39").Value = "The Ice Plane"
$ws.Range("I39").Value = 1856
$ws.Range("J39").Value = 736
